# Apply updated odds values to Sheet1, rows 4-7, as described by the diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 4 ---
$ws.Range("F4").Value  = 3.2
$ws.Range("G4").Value  = 3.45
$ws.Range("H4").Value  = 2.82
$ws.Range("I4").Value  = 3.05
$ws.Range("J4").Value  = 2.78
$ws.Range("K4").Value  = 2.94
$ws.Range("M4").Value  = 1.21
$ws.Range("N4").Value  = 2.06
$ws.Range("O4").Value  = 1.91
$ws.Range("P4").Value  = 1.29
$ws.Range("Q4").Value  = 3.85
$ws.Range("S4").Value  = 8.6
$ws.Range("T4").Value  = 2.7
$ws.Range("U4").Value  = 1.46
$ws.Range("V4").Value  = 1.5
$ws.Range("W4").Value  = 1.4
$ws.Range("Y4").Value  = 6.4
$ws.Range("Z4").Value  = 21
$ws.Range("AD4").Value = 1000
$ws.Range("AG4").Value = 980
$ws.Range("AH4").Value = 980
$ws.Range("AJ4").Value = 980

# --- Row 5 ---
$ws.Range("F5").Value  = 2.1
$ws.Range("K5").Value  = 3.7
$ws.Range("Y5").Value  = 980
$ws.Range("Z5").Value  = 980
$ws.Range("AB5").Value = 8.2
$ws.Range("AF5").Value = 980
$ws.Range("AG5").Value = 980
$ws.Range("AH5").Value = 980
$ws.Range("AJ5").Value = 980
$ws.Range("AK5").Value = 980
$ws.Range("AN5").Value = 980

# --- Row 6 ---
$ws.Range("G6").Value  = 2.14
$ws.Range("I6").Value  = 5.8
$ws.Range("J6").Value  = 2.74
$ws.Range("K6").Value  = 3.4
$ws.Range("M6").Value  = 1.1
$ws.Range("Q6").Value  = 2.24
$ws.Range("W6").Value  = 1.87
$ws.Range("X6").Value  = 980
$ws.Range("Z6").Value  = 980
$ws.Range("AB6").Value = 8.2
$ws.Range("AD6").Value = 980
$ws.Range("AF6").Value = 980
$ws.Range("AG6").Value = 980
$ws.Range("AH6").Value = 980
$ws.Range("AJ6").Value = 980
$ws.Range("AK6").Value = 980
$ws.Range("AN6").Value = 980

# --- Row 7 ---
$ws.Range("F7").Value  = 2.3
$ws.Range("G7").Value  = 2.44
$ws.Range("H7").Value  = 3.5
$ws.Range("J7").Value  = 3.1
$ws.Range("K7").Value  = 3.35
$ws.Range("M7").Value  = 1.13
$ws.Range("N7").Value  = 2.5
$ws.Range("O7").Value  = 1.57
$ws.Range("Q7").Value  = 2.72
$ws.Range("W7").Value  = 1.69
$ws.Range("X7").Value  = 10
$ws.Range("Y7").Value  = 10
$ws.Range("Z7").Value  = 29
$ws.Range("AA7").Value = 120
$ws.Range("AB7").Value = 8.4
$ws.Range("AC7").Value = 8.8
$ws.Range("AD7").Value = 980
$ws.Range("AF7").Value = 980
$ws.Range("AG7").Value = 980
$ws.Range("AH7").Value = 980
$ws.Range("AI7").Value = 120
$ws.Range("AJ7").Value = 980
$ws.Range("AK7").Value = 980
$ws.Range("AN7").Value = 44
$ws.Range("AO7").Value = 140
